# Challenge import: add "Show Statistics Continuously" and "Gameweek" columns
# to the Challenges sheet (Gameweeks import feature).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New header cells
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New data cells for the sample row.
# S2 must be stored as literal text "true" (not a native boolean), so build
# it as a formula result and paste back as a value - this avoids Excel's
# "true"/"false" auto-boolean literal detection while leaving no leftover
# formula or formatting behind.
$ws.Range("ZZ1").Formula = "=""tr""&""ue"""
$ws.Range("ZZ1").Copy()
$ws.Range("S2").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("T2").Value = 1
